$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("B2").Value = 0.0
$ws.Range("C2").Value = 61.50000000000001
$ws.Range("D2").Value = 0.0

# Remove row 3 entirely
$ws.Rows.Item(3).Delete()
